# Business analysis updated (name adjustment):
# Rename the product "Emotion Detector" -> "ShireEye" throughout the
# "Risikoanalyse" worksheet's risk descriptions, and leave the selection
# on cell E4 (where the author ended up after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Risikoanalyse")

$ws.Range("E4").Value = "ShireEye stuft Alter/Emotion falsch ein"
$ws.Range("D6").Value = "ShireEye fällt bis zu einer Woche aus"
$ws.Range("E6").Value = "Hardware des ShireEye wurde beschädigt"
$ws.Range("E7").Value = "Hardware des ShireEye wurde beschädigt oder hat offene Stellen"
$ws.Range("D9").Value = "ShireEye fällt bis zu einer Woche aus"

$ws.Activate()
$ws.Range("E4").Select()
